# Update the "Introduction " sheet
$wsIntro = $excel.ActiveWorkbook.Worksheets.Item("Introduction ")

# Ref Version: V1.2 -> V0.4
$wsIntro.Range("D7").Value = "V0.4"

# Last update: "30/01/2020" text -> real date 2020-09-02 (serial 44076)
$wsIntro.Range("D9").NumberFormat = "mm-dd-yy"
$wsIntro.Range("D9").Value = (Get-Date -Year 2020 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)

# Update the "Cross review points " sheet
$wsCr = $excel.ActiveWorkbook.Worksheets.Item("Cross review points ")

# Update Status column (H) from Open to Resolved for rows 2-9
$wsCr.Range("H2:H9").Value = "Resolved"

# Add new row 10 with the new open point
$wsCr.Range("A10").NumberFormat = "mm-dd-yy"
$wsCr.Range("A10").Value = (Get-Date -Year 2020 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$wsCr.Range("B10").Value = "T.Sharaby"
$wsCr.Range("C10").Value = "V0.4"
$wsCr.Range("D10").Value = "CYRS"
$wsCr.Range("E10").Value = "ALL"
$wsCr.Range("F10").Value = 'Contradiction between Status table and table of history , the last change change is different "3/2/2020 or 4/2/2020"?'
$wsCr.Range("H10").Value = "Open"

# Add new history row (row 14) on Introduction sheet: version 0.2 entry
$wsIntro.Range("B14").Value = 0.2
$wsIntro.Range("C14").Value = "T.Sharaby"
$wsIntro.Range("E14").NumberFormat = "mm-dd-yy"
$wsIntro.Range("E14").Value = (Get-Date -Year 2020 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$wsIntro.Range("G14").Value = "Add one point and update the status for each req "
